$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.763.77"
$ws.Range("E2").Value = "  +1.11%  "
$ws.Range("D3").Value = "1.631.35"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  -0.90%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.49"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.993"
$ws.Range("E7").Value = "  -0.93%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.33"
$ws.Range("E8").Value = "  +1.59%  "
$ws.Range("E9").Value = "  -2.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0611"
$ws.Range("E10").Value = "  -0.08%  "
$ws.Range("E11").Value = "  +1.90%  "
$ws.Range("D12").Value = "1.862.68"
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("D13").Value = "1.629.68"
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("E15").Value = "  +2.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.25"
$ws.Range("E16").Value = "  +0.55%  "
$ws.Range("D17").Value = "27.760.82"
$ws.Range("E17").Value = "  +1.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "233.30"
$ws.Range("E18").Value = "  +2.07%  "
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("E20").Value = "  +0.61%  "
$ws.Range("E21").Value = "  -0.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.39"
$ws.Range("E22").Value = "  -2.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.34"
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("E24").Value = "  -2.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.25"
$ws.Range("E25").Value = "  +1.06%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.58"
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.994"
$ws.Range("E29").Value = "  -0.85%  "
$ws.Range("E30").Value = "  -0.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0482"
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("E32").Value = "  +2.20%  "
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("D34").Value = "1.405.12"
$ws.Range("E34").Value = "  -3.97%  "
$ws.Range("E35").Value = "  +2.17%  "
$ws.Range("E36").Value = "  +0.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0169"
$ws.Range("E37").Value = "  +0.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.876"
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.557"
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.909"
$ws.Range("E40").Value = "  -1.56%  "
$ws.Range("E41").Value = "  +0.29%  "
$ws.Range("E42").Value = "  -0.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "66.63"
$ws.Range("E43").Value = "  -1.56%  "
$ws.Range("E44").Value = "  +6.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.47"
$ws.Range("E45").Value = "  +1.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.20"
$ws.Range("E46").Value = "  -0.46%  "
$ws.Range("D47").Value = "1.771.96"
$ws.Range("E47").Value = "  +0.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.36"
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("E49").Value = "  +0.23%  "
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.58"
$ws.Range("E51").Value = "  -1.24%  "
